# Automatische test-sync: 2025-08-14 20:27:50
# Append a new log row to the "Logs" sheet, extend the conditional
# formatting ranges that cover the data rows, and refresh the summary
# count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Append the new row to the "Logs" sheet ---------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 6
$logs.Cells.Item($newRow, 1).Value  = "Demo inplannen"
$logs.Cells.Item($newRow, 2).Value  = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item($newRow, 3).Value  = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item($newRow, 4).Value  = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item($newRow, 6).Value  = "2025-08-14 20:27:00"
$logs.Cells.Item($newRow, 7).Value  = "Nee"
$logs.Cells.Item($newRow, 8).Value  = "Ja"
$logs.Cells.Item($newRow, 9).Value  = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- 2. Extend the conditional formatting ranges from row 5 to row 6 -----
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "5")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "6")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- 3. Update the summary count on the "Dashboard" sheet ----------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 5
